# "Generate Report for Archive"
#
# The localization-status report is re-generated: the row that used to
# describe 942fa5b4-8a8b-41cb-993c-03167ab1282a (which was "Ready for
# handoff") and the row describing eac70f27-66a6-43f6-b667-1d9cfaee90a4
# (which was "In Translation") swap places - 942fa5b4 now sorts ahead of
# eac70f27, and 942fa5b4's status is refreshed to "In Translation" to
# match the regenerated report, keeping each entry's own handoff/target
# file names and datetimes.
#
# This script applies that change on the "Overview", "zh-cn" and
# "de-de" worksheets, for both the plain cell values and the hyperlink
# display text shown for the markdown / xlf file name columns.

$wb = $excel.ActiveWorkbook

function Set-RowValues($ws, $row, $values) {
    $col = 1
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

function Update-Hyperlink($ws, $row, $col, $newText) {
    foreach ($h in $ws.Hyperlinks) {
        if (($h.Range.Row -eq $row) -and ($h.Range.Column -eq $col)) {
            $h.TextToDisplay = $newText
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: columns A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-RowValues $wsOverview 3 @("942fa5b4-8a8b-41cb-993c-03167ab1282a.md", "In Translation", "In Translation", "2016-25-09 10:25:54")
Set-RowValues $wsOverview 4 @("eac70f27-66a6-43f6-b667-1d9cfaee90a4.md", "In Translation", "In Translation", "2016-23-09 10:23:18")

Update-Hyperlink $wsOverview 3 1 "942fa5b4-8a8b-41cb-993c-03167ab1282a.md"
Update-Hyperlink $wsOverview 4 1 "eac70f27-66a6-43f6-b667-1d9cfaee90a4.md"

# ---------------------------------------------------------------------
# zh-cn sheet: A=Source File Name, B=File Extension, C=Status,
#              D=Latest Handoff File, E=Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-RowValues $wsZhCn 3 @("942fa5b4-8a8b-41cb-993c-03167ab1282a.md", ".md", "In Translation", "942fa5b4-8a8b-41cb-993c-03167ab1282a.d2a9adcaee6d1c80e507967e457b0a6766e83171.zh-cn.xlf", "2016-03-09 10:25:43")
Set-RowValues $wsZhCn 4 @("eac70f27-66a6-43f6-b667-1d9cfaee90a4.md", ".md", "In Translation", "eac70f27-66a6-43f6-b667-1d9cfaee90a4.6095f2be09a451b85982c0b18398bfb76e01cc54.zh-cn.xlf", "2016-03-09 10:22:59")

Update-Hyperlink $wsZhCn 3 1 "942fa5b4-8a8b-41cb-993c-03167ab1282a.md"
Update-Hyperlink $wsZhCn 3 4 "942fa5b4-8a8b-41cb-993c-03167ab1282a.d2a9adcaee6d1c80e507967e457b0a6766e83171.zh-cn.xlf"
Update-Hyperlink $wsZhCn 4 1 "eac70f27-66a6-43f6-b667-1d9cfaee90a4.md"
Update-Hyperlink $wsZhCn 4 4 "eac70f27-66a6-43f6-b667-1d9cfaee90a4.6095f2be09a451b85982c0b18398bfb76e01cc54.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet: A=Source File Name, B=File Extension, C=Status,
#              D=Latest Handoff File, E=Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-RowValues $wsDeDe 3 @("942fa5b4-8a8b-41cb-993c-03167ab1282a.md", ".md", "In Translation", "942fa5b4-8a8b-41cb-993c-03167ab1282a.d2a9adcaee6d1c80e507967e457b0a6766e83171.de-de.xlf", "2016-03-09 10:25:54")
Set-RowValues $wsDeDe 4 @("eac70f27-66a6-43f6-b667-1d9cfaee90a4.md", ".md", "In Translation", "eac70f27-66a6-43f6-b667-1d9cfaee90a4.6095f2be09a451b85982c0b18398bfb76e01cc54.de-de.xlf", "2016-03-09 10:23:18")

Update-Hyperlink $wsDeDe 3 1 "942fa5b4-8a8b-41cb-993c-03167ab1282a.md"
Update-Hyperlink $wsDeDe 3 4 "942fa5b4-8a8b-41cb-993c-03167ab1282a.d2a9adcaee6d1c80e507967e457b0a6766e83171.de-de.xlf"
Update-Hyperlink $wsDeDe 4 1 "eac70f27-66a6-43f6-b667-1d9cfaee90a4.md"
Update-Hyperlink $wsDeDe 4 4 "eac70f27-66a6-43f6-b667-1d9cfaee90a4.6095f2be09a451b85982c0b18398bfb76e01cc54.de-de.xlf"
